# Update the "date" footer placeholder text across the slide master,
# every slide layout and the notes master from 06/02/2023 to 09/03/2023.
# (ppPlaceholderDate = 16)

$p = $ppt.ActivePresentation
$newDate = "09/03/2023"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster

# Every slide layout attached to the slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li)
}

# Notes master
Update-DatePlaceholder $p.NotesMaster
